$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A46").Value = "14/05/2020"
$ws.Range("B46").Value = 1022
$ws.Range("C46").Value = 176

$table = $ws.ListObjects("Condicion_Pacientes")
$table.Resize($ws.Range("A1:F46"))

$ws.Range("A46").HorizontalAlignment = -4152

[void]$ws.Range("C47").Select()
